# Practice tasks and final revisions
# Rename sheets and update task-order file-name cells across all 5 task-order sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "GNG_TO-16509961770360882"
$ws.Range("B2").Value = "go_stims-16509961769960325.csv"
$ws.Range("B3").Value = "GNG_stims-1650996177019994.csv"
$ws.Range("B4").Value = "go_stims-1650996177019994.csv"
$ws.Range("B5").Value = "GNG_stims-16509961770360882.csv"

# --- Sheet 2: NB ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-16509961790687442"
$ws.Range("B2").Value = "ZB-match_0-1650996177684742.csv"
$ws.Range("B3").Value = "TB-16509961790527065.csv"
$ws.Range("B4").Value = "TB-16509961790287454.csv"
$ws.Range("B5").Value = "OB-16509961786287441.csv"
$ws.Range("B6").Value = "OB-16509961780927422.csv"
$ws.Range("B7").Value = "TB-16509961786847062.csv"
$ws.Range("B8").Value = "ZB-match_7-16509961772767105.csv"
$ws.Range("B9").Value = "OB-16509961783407044.csv"
$ws.Range("B10").Value = "ZB-match_7-1650996177788743.csv"

# --- Sheet 3: RS (sheet name only, no cell changes) ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "RS_TO-16509961790687442"

# --- Sheet 4: TOL ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "TOL_TO-16509961791327393"
$ws.Range("B2").Value = "MM_stims-1650996179100739.csv"
$ws.Range("B3").Value = "ZM_stims-16509961790767055.csv"
$ws.Range("B4").Value = "MM_stims-16509961791167068.csv"
$ws.Range("B5").Value = "ZM_stims-1650996179100739.csv"
$ws.Range("B6").Value = "MM_stims-16509961791327393.csv"
$ws.Range("B7").Value = "ZM_stims-16509961791167068.csv"

# --- Sheet 5: vSAT ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "vSAT_TO-16509961792127416"
$ws.Range("B2").Value = "vSAT_stims-16509961791967418.csv"
$ws.Range("B3").Value = "SAT_stims-16509961791327393.csv"
$ws.Range("B4").Value = "SAT_stims-16509961791647093.csv"
$ws.Range("B5").Value = "vSAT_stims-16509961791807442.csv"
